# 🔄 Actualización automática del tracker
# Fill in resultado (G) / profit (H) for newly-settled rows, and fix
# a couple of event_id (A) cells that had been written as text instead
# of numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> @(resultado, profit)
$updates = @{
    102 = @("Fallo",   -1)
    112 = @("Acierto",  1.3)
    113 = @("Acierto",  1.38)
    114 = @("Fallo",   -1)
    115 = @("Fallo",   -1)
    122 = @("Fallo",   -1)
    123 = @("Acierto",  2.5)
    128 = @("Acierto",  1.3)
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    $ws.Range("G$row").Value = $values[0]
    $ws.Range("H$row").Value = $values[1]
}

# event_id values that were stored as text -> convert to real numbers
$ws.Range("A132").Value = 14386752
$ws.Range("A133").Value = 14310235
